$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "57.654.94"
$ws.Range("E2").Value = "  -4.01%  "

# Row 3
$ws.Range("D3").Value = "2.925.53"
$ws.Range("E3").Value = "  -2.25%  "

# Row 4
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").Value = "549.07"
$ws.Range("E5").Value = "  -3.63%  "

# Row 6
$ws.Range("D6").Value = "129.97"
$ws.Range("E6").Value = "  +4.11%  "

# Row 7
$ws.Range("E7").Value = "  -0.08%  "

# Row 8
$ws.Range("D8").Value = "0.511"
$ws.Range("E8").Value = "  +1.88%  "

# Row 9
$ws.Range("D9").Value = "2.917.67"
$ws.Range("E9").Value = "  -2.35%  "

# Row 10
$ws.Range("E10").Value = "  -3.68%  "

# Row 11
$ws.Range("D11").Value = "4.76"
$ws.Range("E11").Value = "  -5.14%  "

# Row 12
$ws.Range("D12").Value = "0.445"
$ws.Range("E12").Value = "  +0.68%  "

# Row 13
$ws.Range("E13").Value = "  -0.51%  "

# Row 14
$ws.Range("D14").Value = "32.72"
$ws.Range("E14").Value = "  +1.14%  "

# Row 15
$ws.Range("E15").Value = "  +1.77%  "

# Row 16
$ws.Range("D16").Value = "3.407.15"
$ws.Range("E16").Value = "  -2.37%  "

# Row 17
$ws.Range("D17").Value = "6.91"
$ws.Range("E17").Value = "  +6.28%  "

# Row 18
$ws.Range("D18").Value = "2.917.29"
$ws.Range("E18").Value = "  -2.42%  "

# Row 19
$ws.Range("D19").Value = "57.608.73"
$ws.Range("E19").Value = "  -4.08%  "

# Row 20
$ws.Range("D20").Value = "417.28"
$ws.Range("E20").Value = "  -2.28%  "

# Row 21
$ws.Range("D21").Value = "13.26"
$ws.Range("E21").Value = "  +1.01%  "

# Row 22
$ws.Range("D22").Value = "0.690"
$ws.Range("E22").Value = "  +3.09%  "

# Row 23
$ws.Range("B23").Value = "InternetComputer(DFINITY)"
$ws.Range("C23").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D23").Value = "13.25"
$ws.Range("E23").Value = "  +2.98%  "

# Row 24
$ws.Range("B24").Value = "Uniswap"
$ws.Range("C24").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D24").Value = "6.98"
$ws.Range("E24").Value = "  -1.00%  "

# Row 25
$ws.Range("D25").Value = "79.67"
$ws.Range("E25").Value = "  +0.36%  "

# Row 26
$ws.Range("E26").Value = "  -0.01%  "

# Row 27
$ws.Range("E27").Value = "  -0.01%  "

# Row 28
$ws.Range("D28").Value = "2.46"
$ws.Range("E28").Value = "  -2.37%  "

# Row 29
$ws.Range("E29").Value = "  +3.00%  "

# Row 30
$ws.Range("D30").Value = "7.39"
$ws.Range("E30").Value = "  +2.97%  "

# Row 31
$ws.Range("D31").Value = "25.24"
$ws.Range("E31").Value = "  +0.17%  "

# Row 32
$ws.Range("D32").Value = "5.98"
$ws.Range("E32").Value = "  -1.72%  "

# Row 33
$ws.Range("D33").Value = "0.0967"
$ws.Range("E33").Value = "  +0.89%  "

# Row 34
$ws.Range("D34").Value = "5.68"
$ws.Range("E34").Value = "  +1.84%  "

# Row 35
$ws.Range("D35").Value = "0.933"
$ws.Range("E35").Value = "  +0.77%  "

# Row 36
$ws.Range("D36").Value = "2.06"
$ws.Range("E36").Value = "  +3.27%  "

# Row 37
$ws.Range("E37").Value = "  -3.99%  "

# Row 38
$ws.Range("D38").Value = "8.73"
$ws.Range("E38").Value = "  +2.83%  "

# Row 39
$ws.Range("D39").Value = "0.0₃0681"
$ws.Range("E39").Value = "  +4.26%  "

# Row 40
$ws.Range("D40").Value = "2.56"
$ws.Range("E40").Value = "  +5.44%  "

# Row 41
$ws.Range("E41").Value = "  -0.24%  "

# Row 42
$ws.Range("D42").Value = "374.59"
$ws.Range("E42").Value = "  +1.02%  "

# Row 43
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "0.0344"
$ws.Range("E43").Value = "  -2.95%  "

# Row 44
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "2.700.84"
$ws.Range("E44").Value = "  +1.31%  "

# Row 46
$ws.Range("D46").Value = "123.75"
$ws.Range("E46").Value = "  +2.48%  "

# Row 47
$ws.Range("E47").Value = "  +1.11%  "

# Row 48
$ws.Range("E48").Value = "  +0.65%  "

# Row 49
$ws.Range("D49").Value = "1.95"
$ws.Range("E49").Value = "  -0.99%  "

# Row 50
$ws.Range("D50").Value = "22.87"
$ws.Range("E50").Value = "  -1.50%  "

# Row 51
$ws.Range("D51").Value = "1.99"
$ws.Range("E51").Value = "  -0.15%  "
